$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.507.59"
$ws.Range("E2").Value = "  -2.97%  "

$ws.Range("D3").Value = "3.487.93"
$ws.Range("E3").Value = "  -0.22%  "

$ws.Range("D4").Value = "'1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.18%  "

$ws.Range("D5").Value = "'554.47"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.87%  "

$ws.Range("D6").Value = "'178.53"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -5.11%  "

$ws.Range("D7").Value = "'0.639"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +4.10%  "

$ws.Range("E8").Value = "  +0.05%  "

$ws.Range("E9").Value = "  -1.12%  "

$ws.Range("E10").Value = "  +2.65%  "

$ws.Range("D11").Value = "'53.68"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -5.12%  "

$ws.Range("E12").Value = "  -2.02%  "

$ws.Range("D13").Value = "'9.25"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.54%  "

$ws.Range("D14").Value = "4.041.11"
$ws.Range("E14").Value = "  +0.06%  "

$ws.Range("D15").Value = "'18.69"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.99%  "

$ws.Range("D16").Value = "3.487.03"
$ws.Range("E16").Value = "  +0.04%  "

$ws.Range("E17").Value = "  +0.30%  "

$ws.Range("D18").Value = "'12.11"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.28%  "

$ws.Range("D19").Value = "65.507.05"
$ws.Range("E19").Value = "  -3.20%  "

$ws.Range("D20").Value = "'0.990"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.14%  "

$ws.Range("D21").Value = "'417.49"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.90%  "

$ws.Range("E22").Value = "  +2.78%  "

$ws.Range("D23").Value = "'86.21"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.56%  "

$ws.Range("D24").Value = "'4.12"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.61%  "

$ws.Range("D25").Value = "'12.87"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +8.76%  "

$ws.Range("D26").Value = "'10.80"
$ws.Range("D26").Style = "Normal"

$ws.Range("E27").Value = "  -2.99%  "

$ws.Range("D28").Value = "'6.03"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.89%  "

$ws.Range("D29").Value = "'9.05"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +4.22%  "

$ws.Range("D30").Value = "'30.24"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.33%  "

$ws.Range("D31").Value = "'6.50"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -5.51%  "

$ws.Range("D32").Value = "'611.48"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -9.68%  "

$ws.Range("D33").Value = "'11.73"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.15%  "

$ws.Range("E34").Value = "  -0.82%  "

$ws.Range("D35").Value = "'59.45"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.56%  "

$ws.Range("D36").Value = "'0.147"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +9.66%  "

$ws.Range("D37").Value = "'0.999"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.09%  "

$ws.Range("D38").Value = "'37.43"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.79%  "

$ws.Range("D39").Value = "3.393.03"
$ws.Range("E39").Value = "  +10.97%  "

$ws.Range("E40").Value = "  -5.69%  "

$ws.Range("E41").Value = "  -5.91%  "

$ws.Range("D42").Value = "'0.999"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.03%  "

$ws.Range("D43").Value = "'3.23"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -7.34%  "

$ws.Range("D44").Value = "'2.82"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -5.15%  "

$ws.Range("E45").Value = "  -9.41%  "

$ws.Range("B46").Value = "VeChain"
$ws.Range("C46").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D46").Value = "'0.0414"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.79%  "

$ws.Range("B47").Value = "ApeXProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D47").Value = "'3.23"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.36%  "

$ws.Range("D48").Value = "'2.70"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.39%  "

$ws.Range("D49").Value = "'0.133"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.02%  "

$ws.Range("D50").Value = "'8.50"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.24%  "

$ws.Range("D51").Value = "'137.61"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.56%  "
